# edit.ps1 - "edit survival figure - font size and axis"
#
# Extends the underlying data table from columns D:M (10 days) out to
# D:U (18 days) and updates the mortality chart (category + value series
# ranges, several font sizes, title/legend layout, and the anchor/size of
# the chart on the sheet) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Worksheet data: extend the "Day" header row (row 2) from column M out to
#    column U (days 10-17), and backfill the three mortality-count rows
#    (3, 4, 5) with zeros for the new columns.
# ---------------------------------------------------------------------------

$dayCols = @("N", "O", "P", "Q", "R", "S", "T", "U")
$dayVal = 10
foreach ($col in $dayCols) {
    $cell = $ws.Range($col + "2")
    $cell.Value = $dayVal
    $cell.NumberFormat = "General"
    $dayVal++
}

foreach ($row in 3, 4, 5) {
    foreach ($col in $dayCols) {
        $ws.Range($col + [string]$row).Value = 0
    }
}

# ---------------------------------------------------------------------------
# 2. Sheet view: scroll down a bit and move the active selection, matching
#    the saved view state in the edited workbook.
# ---------------------------------------------------------------------------

$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("R9").Select()

# ---------------------------------------------------------------------------
# 3. Chart: widen every series' category + value ranges from D:M to D:U.
# ---------------------------------------------------------------------------

$chartObj = $ws.ChartObjects().Item(1)
$chart = $chartObj.Chart

$chart.SeriesCollection().Item(1).Formula = "=SERIES(Sheet1!`$C`$3,Sheet1!`$D`$2:`$U`$2,Sheet1!`$D`$3:`$U`$3,1)"
$chart.SeriesCollection().Item(2).Formula = "=SERIES(Sheet1!`$C`$4,Sheet1!`$D`$2:`$U`$2,Sheet1!`$D`$4:`$U`$4,2)"
$chart.SeriesCollection().Item(3).Formula = "=SERIES(Sheet1!`$C`$5,Sheet1!`$D`$2:`$U`$2,Sheet1!`$D`$5:`$U`$5,3)"

# ---------------------------------------------------------------------------
# 4. Grow the chart's frame on the sheet so the wider plot still fits
#    (bottom-right corner only -- the top-left anchor is unchanged).
# ---------------------------------------------------------------------------

$chartObj.Width = 564.2431640625
$chartObj.Height = 359

# ---------------------------------------------------------------------------
# 5. Chart text sizes: bump up the category axis title ("Day"), the value
#    axis title ("Number of mortalities"), both axes' tick-label text, and
#    the legend text.
# ---------------------------------------------------------------------------

$catAxis = $chart.Axes(1)      # xlCategory
$valAxis = $chart.Axes(2)      # xlValue

$catAxis.AxisTitle.Font.Size = 20
$catAxis.TickLabels.Font.Size = 16

$valAxis.AxisTitle.Font.Size = 20
$valAxis.TickLabels.Font.Size = 16

$chart.Legend.Font.Size = 18

# Re-position the "Day" category-axis title and "Number of mortalities"
# value-axis title using the manual layout coordinates saved with the
# resized chart.
$catAxis.AxisTitle.Left = 0.46974418927971084 * $chartObj.Width + $chartObj.Left
$catAxis.AxisTitle.Top = 0.91823196543329033 * $chartObj.Height + $chartObj.Top

$valAxis.AxisTitle.Left = 0.0029413177285423598 * $chartObj.Width + $chartObj.Left
$valAxis.AxisTitle.Top = 0.28299585462680682 * $chartObj.Height + $chartObj.Top

